# "few changes and scrum details in case study folder"
#
# The Sprint 1 retrospective sheet had its "Good / Bad / Ugly" content
# cells filled in with the actual retro notes (they were empty
# placeholders before). Fill the four comment cells and leave the
# selection on the last-edited cell (C11), matching the author's save
# state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Good -> Continue Doing
$ws.Range("C5").Value = "the coding is done at good pace and keep up"

# The Bad -> Stop Doing
$ws.Range("C10").Value = "versioning isssues with git"

# The Ugly -> Start Doing
$ws.Range("C16").Value = "stop pushing to master directly "

# The Bad -> Stop Doing (second bullet)
$ws.Range("C11").Value = "not using branches for versioning"

# Leave the selection where the author left it when saving.
[void]$ws.Range("C11:L11").Select()
